{"js": "// Change the \"Last modified on\" date from 10/31/2018 to 2/27/2019,\n// and rename the \"AUS-DB-STG-03_Local Administrators\" group to\n// \"AUS-DB-STG-03_ Administrators\" (dropping \"Local\"), matching the\n// author's edit described in the commit message / diff.\n\nconst body = context.document.body;\n\n// --- 1. Update the \"Last modified on ...\" date line -----------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet dateParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Last modified on\") === 0) {\n    dateParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (dateParagraph) {\n  // Replacing the whole paragraph range (not just a sub-match) also\n  // clears out the old (now stale) \"_GoBack\" bookmark that previously\n  // sat inside this paragraph, matching the target document where that\n  // bookmark is gone from here entirely.\n  const dateRange = dateParagraph.getRange();\n  dateRange.insertText(\"Last modified on 2/27/2019\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. Rename the administrators group text -------------------------\nconst groupMatches = body.search(\"AUS-DB-STG-03_Local Administrators\", { matchCase: true });\nawait context.sync();\n\nif (groupMatches.items.length > 0) {\n  const fullMatch = groupMatches.items[0];\n\n  // Narrow the search to just the word \"Local\" within the match so we\n  // keep the surrounding space (\u2026STG-03_ Administrators, with a single\n  // space where \"Local\" used to be).\n  const localMatches = fullMatch.search(\"Local\", { matchCase: true });\n  await context.sync();\n\n  if (localMatches.items.length > 0) {\n    const localRange = localMatches.items[0];\n    localRange.insertText(\"\", Word.InsertLocation.replace);\n    // Word leaves its \"_GoBack\" bookmark (last-edit marker) at the spot\n    // of the most recent edit -- recreate that here now that \"Local\"\n    // has been removed.\n    localRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Change the \"Last modified on\" date from 10/31/2018 to 2/27/2019,\n# and rename the \"AUS-DB-STG-03_Local Administrators\" group to\n# \"AUS-DB-STG-03_ Administrators\" (dropping \"Local\"), matching the\n# author's edit described in the commit message / diff.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the \"Last modified on ...\" date line ---------------------\n$dateRange = $d.Content\n$dateFind = $dateRange.Find\n$dateFind.Text = \"Last modified on 10/31/2018\"\n[void]$dateFind.Execute()\n\nif ($dateFind.Found) {\n    # Replacing the whole matched range's text (instead of a narrower\n    # sub-match) also clears out the old, now-stale \"_GoBack\" bookmark\n    # that previously sat inside this paragraph, matching the target\n    # document where that bookmark is gone from here entirely.\n    $dateRange.Text = \"Last modified on 2/27/2019\"\n}\n\n# --- 2. Rename the administrators group text -----------------------------\n$groupRange = $d.Content\n$groupFind = $groupRange.Find\n$groupFind.Text = \"AUS-DB-STG-03_Local Administrators\"\n[void]$groupFind.Execute()\n\nif ($groupFind.Found) {\n    # Narrow to just the word \"Local\" (on a duplicate range, since reusing\n    # the same Range/Find object searches forward from the current match\n    # rather than within it) so the surrounding space is preserved,\n    # producing \"...STG-03_ Administrators\".\n    $localRange = $groupRange.Duplicate\n    $localFind = $localRange.Find\n    $localFind.Text = \"Local\"\n    [void]$localFind.Execute()\n\n    if ($localFind.Found) {\n        $localRange.Text = \"\"\n\n        # Word leaves its \"_GoBack\" bookmark (last-edit marker) at the spot\n        # of the most recent edit -- recreate that here now that \"Local\"\n        # has been removed (localRange is now collapsed at that point).\n        [void]$d.Bookmarks.Add(\"_GoBack\", $localRange)\n    }\n}\n"}
